$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 212, pushing existing rows 212..283 down to 213..284
$ws.Rows(212).Insert()

$ws.Range("A212").Value = 3
$ws.Range("B212").Value = "Femacal de La Calera"
$ws.Range("C212").Value = "Coquimbo"
$ws.Range("D212").Value = 44588
$ws.Range("E212").Value = 5
$ws.Range("F212").Value = 100112043
$ws.Range("G212").Value = "Pepino ensalada"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 90
$ws.Range("K212").Value = 11000
$ws.Range("L212").Value = 12000
$ws.Range("M212").Value = 11444
$ws.Range("N212").Value = "$/caja 70 unidades"
$ws.Range("O212").Value = "Limache"
$ws.Range("P212").Value = 163
$ws.Range("Q212").Value = 70
$ws.Range("R212").Value = "Hortaliza"
